$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("D5").Value = 44602
$ws.Range("J5").Value = 240
$ws.Range("K5").Value = 23000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 23500
$ws.Range("P5").Value = 940

# Row 6 updates
$ws.Range("D6").Value = 44293
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 23000
$ws.Range("M6").Value = 21500
$ws.Range("N6").Value = "$/malla 25 kilos"
$ws.Range("P6").Value = 860
$ws.Range("Q6").Value = 25
